$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.843.32'
$ws.Range("E2").Value = '  +3.86%  '

$ws.Range("D3").Value = '3.020.00'
$ws.Range("E3").Value = '  +3.03%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '565.51'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.97%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.77'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +7.95%  '

$ws.Range("E7").Value = '  -0.11%  '

$ws.Range("E8").Value = '  +1.72%  '

$ws.Range("D9").Value = '3.010.32'
$ws.Range("E9").Value = '  +2.89%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.134'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.27%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.31'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +10.94%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.460'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.44%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000232'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.28%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.22'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.00%  '

$ws.Range("E15").Value = '  +1.79%  '

$ws.Range("D16").Value = '3.523.69'
$ws.Range("E16").Value = '  +3.05%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.21'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +5.45%  '

$ws.Range("D18").Value = '3.020.25'
$ws.Range("E18").Value = '  +3.04%  '

$ws.Range("D19").Value = '59.834.47'
$ws.Range("E19").Value = '  +3.80%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '437.67'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.05%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.72'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.18%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.722'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.51%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.14'
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.34'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.47%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '80.82'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.27%  '

$ws.Range("E26").Value = '  -0.05%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.26'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +13.56%  '

$ws.Range("E28").Value = '  +0.09%  '

$ws.Range("E29").Value = '  +3.54%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.85'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.04%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.08'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.87%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.29'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.07%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.103'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.58%  '

$ws.Range("D34").Value = '0.0₃0790'
$ws.Range("E34").Value = '  +16.09%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.01'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +7.64%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.93'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.84%  '

$ws.Range("E37").Value = '  +2.56%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '49.18'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.59%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.67'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.32%  '

$ws.Range("E40").Value = '  +10.11%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '405.80'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +7.71%  '

$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '2.786.45'
$ws.Range("E42").Value = '  +3.84%  '

$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0354'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.54%  '

$ws.Range("E44").Value = '  -0.35%  '

$ws.Range("E45").Value = '  +6.46%  '

$ws.Range("E46").Value = '  +0.02%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '122.99'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.81%  '

$ws.Range("E48").Value = '  +3.30%  '

$ws.Range("E49").Value = '  +1.53%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '33.92'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +19.91%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.66'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.62%  '
